$d = $word.ActiveDocument

# --- Locate the "Onsdag 22 Februari" paragraph and the one that follows it
#     (the IT-JURIDIK paragraph). The document, as authored, currently ends
#     right after the IT-JURIDIK paragraph.
$foundOnsdag = $false
$onsdagIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Onsdag 22 Februari*") {
        $onsdagIndex = $i
        $foundOnsdag = $true
    }
}

$itJuridikIndex = $onsdagIndex + 1

# The "_GoBack" bookmark currently sits right at the end of the "Onsdag 22
# Februari" paragraph; it needs to move to the end of the new content being
# added below, so remove it from its old location first.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$pIT = $d.Paragraphs.Item($itJuridikIndex)

# Create a single empty paragraph right after the IT-JURIDIK paragraph; the
# new content will be written into it (and beyond) via InsertXML below.
$pIT.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs.Item($itJuridikIndex + 1)

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:rPr><w:b/></w:rPr></w:pPr>
            <w:r><w:rPr><w:b/></w:rPr><w:t>Torsdag 23 Februari</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:ind w:left="1304" w:firstLine="1"/></w:pPr>
            <w:r><w:t>Blev klara med &#8217;&#228;ndra&#8217; inst&#228;llningar(50 minuter)</w:t></w:r>
            <w:r><w:br/><w:t>Blev kaxiga och gjorde avboka procedur(20 minuter)</w:t></w:r>
            <w:r><w:br/><w:t>Avboknings kaos(</w:t></w:r>
            <w:r><w:t xml:space="preserve">40 </w:t></w:r>
            <w:r><w:t>minuter</w:t></w:r>
            <w:r><w:t>)</w:t></w:r>
            <w:r><w:br/><w:t>Sammanfattning n&#228;stan f&#228;rdig(20 minuter)</w:t></w:r>
            <w:r><w:br/><w:t>Avbokning f&#228;rdig(250 minuter)</w:t></w:r>
            <w:r><w:br/><w:t>Raderade oanv&#228;nd data(10 minuter)</w:t></w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:br/></w:r>
            <w:r><w:br/></w:r>
            <w:r><w:br/></w:r>
            <w:r><w:tab/></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$pNew.Range.InsertXML($xmlSnippet)

# InsertXML leaves one extra (empty, unformatted) trailing paragraph behind
# after the inserted content. Before merging it away by deleting its
# preceding paragraph mark, copy the real content paragraph's indentation
# onto it so that the merge keeps the indentation (Word otherwise adopts the
# formatting of the paragraph whose mark "survives" the merge).
$pContent = $d.Paragraphs.Item($itJuridikIndex + 2)
$pExtra = $d.Paragraphs.Item($itJuridikIndex + 3)
$pExtra.Range.ParagraphFormat.LeftIndent = $pContent.Range.ParagraphFormat.LeftIndent
$pExtra.Range.ParagraphFormat.FirstLineIndent = $pContent.Range.ParagraphFormat.FirstLineIndent

$boundary = $d.Range($pContent.Range.End - 1, $pContent.Range.End)
$boundary.Delete()
